$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 held a stray "984972 - Hugo Ricardo Zschommler Sandim"
# value (no label in column A). Delete that whole row so everything below it
# shifts up by one.
$ws.Rows("13:13").Delete()

# After the shift, a handful of data cells (columns B/C) now hold content
# that no longer matches their row - update them to their final text.

# Row 10 "Objetivos:" -> docente name
$ws.Range("B10").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C10").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 13 "Programa resumido:" -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 "Programa:" -> activation date (copy from the "Ativação:" row so the
# "01/01/2012" text is carried over verbatim instead of being reinterpreted
# as a date serial number by a fresh .Value assignment)
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 "Método:" -> docente name
$ws.Range("B18").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C18").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 19 "Critério:" -> teaching method text
$ws.Range("B19").Value = "Aulas expositivas e aulas práticas de demonstração em oficina. Visita a feiras."
$ws.Range("C19").Value = "Aulas expositivas e aulas práticas de demonstração em oficina. Visita a feiras."

# Row 20 "Norma de recuperação:" -> evaluation criteria text
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# Row 21 "Bibliografia:" -> recovery norm text
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
